$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to remain text even though many values look numeric.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.953.19"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "1.893.71"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "0.7770"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("D6").Value = "243.79"

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "0.3132"
$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").Value = "25.82"
$ws.Range("E9").Value = "  +1.86%  "

$ws.Range("D10").Value = "0.07240"
$ws.Range("E10").Value = "  +0.17%  "

$ws.Range("D11").Value = "0.08714"
$ws.Range("E11").Value = "  +7.86%  "

$ws.Range("D12").Value = "2.098.68"
$ws.Range("E12").Value = "  +9.50%  "

$ws.Range("D13").Value = "0.7743"
$ws.Range("E13").Value = "  +1.15%  "

$ws.Range("D14").Value = "5.420"
$ws.Range("E14").Value = "  -1.26%  "

$ws.Range("D15").Value = "94.49"
$ws.Range("E15").Value = "  +2.32%  "

$ws.Range("D16").Value = "30.346.99"
$ws.Range("E16").Value = "  +1.68%  "

$ws.Range("D17").Value = "6.190"
$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("D18").Value = "2.371.15"
$ws.Range("E18").Value = "  +9.42%  "

$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("D20").Value = "246.05"
$ws.Range("E20").Value = "  +1.04%  "

$ws.Range("D21").Value = "0.000007856"
$ws.Range("E21").Value = "  +1.15%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").Value = "8.118"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").Value = "0.1651"
$ws.Range("E25").Value = "  +5.61%  "

$ws.Range("D26").Value = "9.481"
$ws.Range("E26").Value = "  +0.89%  "

$ws.Range("D27").Value = "163.27"

$ws.Range("E28").Value = "  +0.52%  "

$ws.Range("E29").Value = "  +0.54%  "

$ws.Range("E30").Value = "  -0.32%  "

$ws.Range("E31").Value = "  -0.52%  "

$ws.Range("D32").Value = "4.515"

$ws.Range("D33").Value = "4.141"
$ws.Range("E33").Value = "  +0.98%  "

$ws.Range("D34").Value = "0.05474"
$ws.Range("E34").Value = "  -0.45%  "

$ws.Range("E35").Value = "  -1.18%  "

$ws.Range("E36").Value = "  +0.98%  "

$ws.Range("E37").Value = "  +0.84%  "

$ws.Range("D38").Value = "2.704"
$ws.Range("E38").Value = "  +2.77%  "

$ws.Range("D39").Value = "0.01973"
$ws.Range("E39").Value = "  +2.93%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  +2.20%  "

$ws.Range("D42").Value = "1.113.04"
$ws.Range("E42").Value = "  -2.40%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "6.123"
$ws.Range("E43").Value = "  +4.00%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "73.61"
$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("D45").Value = "2.252.03"
$ws.Range("E45").Value = "  +9.65%  "

$ws.Range("D46").Value = "0.8505"
$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "104.03"
$ws.Range("E47").Value = "  +0.08%  "

$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("D49").Value = "1.880"
$ws.Range("E49").Value = "  -0.32%  "

$ws.Range("D50").Value = "7.628"
$ws.Range("E50").Value = "  +2.31%  "

$ws.Range("D51").Value = "9.905"
$ws.Range("E51").Value = "  +0.18%  "

# Remove the forced text format so column D keeps the workbook default style (no explicit numFmt).
$ws.Range("D2:D51").Style = "Normal"
